$wb = $excel.ActiveWorkbook

# Add the new "dailyQuestStyle" sheet after "dailyQuestStar"
$afterSheet = $wb.Worksheets.Item("dailyQuestStar")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "dailyQuestStyle"

